$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 3694
$ws.Range("J3").Value = 3878
$ws.Range("G4").Value = 1466
$ws.Range("H4").Value = 1695
$ws.Range("J4").Value = 861
$ws.Range("J5").Value = 309
$ws.Range("J6").Value = 4563
$ws.Range("G7").Value = 24690
$ws.Range("H7").Value = 26004
$ws.Range("J7").Value = 13305

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J5").Value = 37
$ws.Range("J8").Value = 863
$ws.Range("J11").Value = 196
$ws.Range("J14").Value = 55
$ws.Range("J15").Value = 152
$ws.Range("J19").Value = 396
$ws.Range("J20").Value = 268
$ws.Range("J21").Value = 28
$ws.Range("J24").Value = 42
$ws.Range("J26").Value = 24
$ws.Range("J29").Value = 755
$ws.Range("J33").Value = 604
$ws.Range("J35").Value = 18
$ws.Range("J37").Value = 417
$ws.Range("J41").Value = 81
$ws.Range("J42").Value = 515
$ws.Range("J43").Value = 120
$ws.Range("J44").Value = 96
$ws.Range("J45").Value = 20
$ws.Range("J46").Value = 48
$ws.Range("J48").Value = 140
$ws.Range("J49").Value = 93
$ws.Range("J50").Value = 76
$ws.Range("J51").Value = 180
$ws.Range("J52").Value = 368
$ws.Range("H53").Value = 246
$ws.Range("J53").Value = 124
$ws.Range("J54").Value = 257
$ws.Range("G63").Value = 223
$ws.Range("J63").Value = 60
$ws.Range("J64").Value = 89
$ws.Range("J65").Value = 349
$ws.Range("J66").Value = 38
$ws.Range("J67").Value = 510
$ws.Range("J72").Value = 54
$ws.Range("J73").Value = 117
$ws.Range("J76").Value = 190
$ws.Range("J78").Value = 176
$ws.Range("J79").Value = 384
$ws.Range("J80").Value = 23
$ws.Range("J85").Value = 607
$ws.Range("J86").Value = 72
$ws.Range("J88").Value = 145
$ws.Range("J89").Value = 171
$ws.Range("J90").Value = 154
$ws.Range("J93").Value = 63
$ws.Range("J95").Value = 208
$ws.Range("J96").Value = 147
$ws.Range("J97").Value = 89
$ws.Range("J99").Value = 192
$ws.Range("G101").Value = 24690
$ws.Range("H101").Value = 26004
$ws.Range("J101").Value = 13305

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 195
$ws.Range("J4").Value = 29
$ws.Range("J6").Value = 195
$ws.Range("J7").Value = 604

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 240
$ws.Range("J6").Value = 182
$ws.Range("J7").Value = 755

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 153
$ws.Range("J4").Value = 49
$ws.Range("J6").Value = 168
$ws.Range("J7").Value = 607

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 79
$ws.Range("J3").Value = 64
$ws.Range("J7").Value = 208

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 112
$ws.Range("J3").Value = 111
$ws.Range("J6").Value = 251
$ws.Range("J7").Value = 515

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 146
$ws.Range("J7").Value = 417

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 45
$ws.Range("J7").Value = 147

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 74
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 49
$ws.Range("J3").Value = 48
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 114
$ws.Range("J3").Value = 208
$ws.Range("J6").Value = 142
$ws.Range("J7").Value = 510

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 48
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 99
$ws.Range("J7").Value = 396

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 104
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 161
$ws.Range("J7").Value = 368

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 273
$ws.Range("J6").Value = 266
$ws.Range("J7").Value = 863

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 349

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 110
$ws.Range("J3").Value = 143
$ws.Range("J6").Value = 102
$ws.Range("J7").Value = 384

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 33
$ws.Range("J3").Value = 38
$ws.Range("J6").Value = 102
$ws.Range("J7").Value = 190

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 16
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J5").Value = 1
$ws.Range("J6").Value = 62
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 46
$ws.Range("J3").Value = 60
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 140

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 120
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 33
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 47
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 56
$ws.Range("J7").Value = 192

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 23

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J3").Value = 22
$ws.Range("H4").Value = 26
$ws.Range("H7").Value = 246
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 35
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 18
Write-Host "Applied 2023-07-03 data update across all sheets."
